$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting/style of the last existing data row (224) down into
# the three new rows (225-227) so the new "A" column cells keep the same
# date style (border/alignment/number format) as the rest of the column.
$ws.Range("A224").Copy()
$ws.Range("A225:A227").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 225
$ws.Range("A225").Value = 45108
$ws.Range("B225").Value = 622009300000
$ws.Range("C225").Value = 0.2230848168473654
$ws.Range("D225").Value = 138760830767.8579

# Row 226
$ws.Range("A226").Value = 45139
$ws.Range("B226").Value = 626029700000
$ws.Range("C226").Value = 0.2195775328268412
$ws.Range("D226").Value = 137462057002.3275

# Row 227
$ws.Range("A227").Value = 45170
$ws.Range("B227").Value = 640762400000
$ws.Range("C227").Value = 0.2126709342634142
$ws.Range("D227").Value = 136271538248.8675
